$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 165
$ws.Range("I2").Value = 165
$ws.Range("K2").Value = 165
$ws.Range("M2").Value = -52

$ws.Range("H17").Value = 1963.3334
$ws.Range("J17").Value = 2056
$ws.Range("L17").Value = 6168
$ws.Range("N17").Value = -6504

$ws.Range("H41").Value = 1954.8334
$ws.Range("I41").Value = 247.75
$ws.Range("K41").Value = 247.75
$ws.Range("M41").Value = 192.25

$ws.Range("H53").Value = 565.58826
$ws.Range("J53").Value = 739.5
$ws.Range("L53").Value = 739.5
$ws.Range("N53").Value = -2013.5

$ws.Range("H62").Value = 4749.5
$ws.Range("J62").Value = 5333
$ws.Range("L62").Value = 5333
$ws.Range("N62").Value = -6581

$ws.Range("H65").Value = 4749.5
$ws.Range("J65").Value = 5333
$ws.Range("L65").Value = 26665
$ws.Range("N65").Value = -32905

$ws.Range("H107").Value = 1261.3125
$ws.Range("I107").Value = 1215
$ws.Range("K107").Value = 1215
$ws.Range("M107").Value = 705

$ws.Range("H111").Value = 7075
$ws.Range("J111").Value = 12396.6
$ws.Range("L111").Value = 37189.8
$ws.Range("N111").Value = -43323.8

$ws.Range("H113").Value = 2566.2727
$ws.Range("I113").Value = 1904.1428
$ws.Range("K113").Value = 1904.1428
$ws.Range("M113").Value = 1349.8572

$ws.Range("H116").Value = 23615706
$ws.Range("J116").Value = 5183.5
$ws.Range("L116").Value = 5183.5
$ws.Range("N116").Value = -12067.5

$ws.Range("H118").Value = 538.25
$ws.Range("I118").Value = 361.5
$ws.Range("K118").Value = 1084.5
$ws.Range("M118").Value = 572.5

$ws.Range("H137").Value = 18072508
$ws.Range("I137").Value = 910800.5600000001
$ws.Range("J137").Value = 41669856
$ws.Range("K137").Value = 2732401.68
$ws.Range("L137").Value = 125009568
$ws.Range("M137").Value = -2729851.68
$ws.Range("N137").Value = -125014668

$ws.Range("H141").Value = 4362.5
$ws.Range("I141").Value = 2128.5715
$ws.Range("K141").Value = 6385.7145
$ws.Range("M141").Value = -1205.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H74").Value = 1342.25
$ws.Range("I74").Value = 1328.9231
$ws.Range("K74").Value = 1328.9231
$ws.Range("M74").Value = -454.9231

$ws.Range("H77").Value = 1342.25
$ws.Range("I77").Value = 1328.9231
$ws.Range("K77").Value = 6644.6155
$ws.Range("M77").Value = -2276.6155

$ws.Range("H122").Value = 4065.5264
$ws.Range("I122").Value = 3471.1428
$ws.Range("K122").Value = 10413.4284
$ws.Range("M122").Value = -7963.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 224988
$ws.Range("J42").Value = 224988
$ws.Range("L42").Value = 224988
$ws.Range("N42").Value = -225644

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1003.875
$ws.Range("I16").Value = 682.46155
$ws.Range("K16").Value = 682.46155
$ws.Range("M16").Value = -395.46155

$ws.Range("H97").Value = 38975
$ws.Range("J97").Value = 38975
$ws.Range("L97").Value = 38975
$ws.Range("N97").Value = -40957

$ws.Range("H104").Value = 59999
$ws.Range("I104").Value = 59999
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 59999
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -57378
$ws.Range("N104").ClearContents()

$ws.Range("H113").Value = 1003.875
$ws.Range("I113").Value = 682.46155
$ws.Range("K113").Value = 682.46155
$ws.Range("M113").Value = 1487.53845

$ws.Range("H115").Value = 49246.668
$ws.Range("J115").Value = 49246.668
$ws.Range("L115").Value = 49246.668
$ws.Range("N115").Value = -51596.668

$ws.Range("H122").Value = 4158.409
$ws.Range("I122").Value = 2666
$ws.Range("K122").Value = 7998
$ws.Range("M122").Value = -5548

$ws.Range("H132").Value = 10758342
$ws.Range("I132").Value = 16668715
$ws.Range("K132").Value = 50006145
$ws.Range("M132").Value = -50003615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20469752
$ws.Range("I4").Value = 68370.125
$ws.Range("K4").Value = 205110.375
$ws.Range("M4").Value = -204998.375

$ws.Range("H5").Value = 796.9375
$ws.Range("I5").Value = 475.8889
$ws.Range("K5").Value = 1427.6667
$ws.Range("M5").Value = -1315.6667

$ws.Range("H34").Value = 1287.8462
$ws.Range("J34").Value = 1941.5714
$ws.Range("L34").Value = 5824.7142
$ws.Range("N34").Value = -5992.7142

$ws.Range("H125").Value = 25000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H130").Value = 370000
$ws.Range("I130").Value = 505000
$ws.Range("K130").Value = 1515000
$ws.Range("M130").Value = -1509980

$ws.Range("H131").Value = 6164869.5
$ws.Range("I131").Value = 91848
$ws.Range("J131").Value = 7316649
$ws.Range("K131").Value = 275544
$ws.Range("L131").Value = 21949947
$ws.Range("M131").Value = -270504
$ws.Range("N131").Value = -21960027

$ws.Range("H135").Value = 796.9375
$ws.Range("I135").Value = 475.8889
$ws.Range("K135").Value = 4283.0001
$ws.Range("M135").Value = -1748.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 17833.334
$ws.Range("J39").Value = 17833.334
$ws.Range("L39").Value = 17833.334
$ws.Range("N39").Value = -18897.334

$ws.Range("H107").Value = 509.46667
$ws.Range("I107").Value = 217.125
$ws.Range("J107").Value = 843.5714
$ws.Range("K107").Value = 217.125
$ws.Range("L107").Value = 843.5714
$ws.Range("M107").Value = 1702.875
$ws.Range("N107").Value = -4683.5714

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H122").Value = 389242.16
$ws.Range("I122").Value = 1113130
$ws.Range("J122").Value = 6007.4116
$ws.Range("K122").Value = 3339390
$ws.Range("L122").Value = 18022.2348
$ws.Range("M122").Value = -3336940
$ws.Range("N122").Value = -22922.2348

$ws.Range("H133").Value = 54772.332
$ws.Range("J133").Value = 54772.332
$ws.Range("L133").Value = 54772.332
$ws.Range("N133").Value = -64892.332

$ws.Range("H134").Value = 550362.7
$ws.Range("J134").Value = 550362.7
$ws.Range("L134").Value = 1651088.1
$ws.Range("N134").Value = -1656158.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4069.3096
$ws.Range("I7").Value = 3542.1694
$ws.Range("J7").Value = 5313.36
$ws.Range("K7").Value = 3542.1694
$ws.Range("L7").Value = 5313.36
$ws.Range("M7").Value = -3430.1694
$ws.Range("N7").Value = -5537.36

$ws.Range("H69").Value = 35000
$ws.Range("J69").Value = 35000
$ws.Range("L69").Value = 35000
$ws.Range("N69").Value = -36622

$ws.Range("H72").Value = 35000
$ws.Range("J72").Value = 35000
$ws.Range("L72").Value = 105000
$ws.Range("N72").Value = -113112

$ws.Range("H122").Value = 6851.9766
$ws.Range("I122").Value = 4110.6924
$ws.Range("K122").Value = 12332.0772
$ws.Range("M122").Value = -9882.0772

$ws.Range("H126").Value = 4069.3096
$ws.Range("I126").Value = 3542.1694
$ws.Range("J126").Value = 5313.36
$ws.Range("K126").Value = 10626.5082
$ws.Range("L126").Value = 15940.08
$ws.Range("M126").Value = -8156.5082
$ws.Range("N126").Value = -20880.08

$ws.Range("H132").Value = 4010.0125
$ws.Range("I132").Value = 3268.6667
$ws.Range("K132").Value = 9806.000100000001
$ws.Range("M132").Value = -7276.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 809.375
$ws.Range("I107").Value = 667.7778
$ws.Range("J107").Value = 991.4286
$ws.Range("K107").Value = 2003.3334
$ws.Range("L107").Value = 2974.2858
$ws.Range("M107").Value = -83.33339999999998
$ws.Range("N107").Value = -6814.2858

$ws.Range("H132").Value = 1599.1882
$ws.Range("I132").Value = 835.5909
$ws.Range("K132").Value = 2506.7727
$ws.Range("M132").Value = 23.22730000000001
